$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows (File_Name, Cell, Num_Inclusions) for rows 2-65,
# reflecting the improved cell-segmentation/histogram-based re-analysis.
$data = @(
    @("1K_PFF_ADAMTS19_01", 0, 18),
    @("1K_PFF_ADAMTS19_01", 1, 2),
    @("1K_PFF_ADAMTS19_02", 0, 1),
    @("1K_PFF_ADAMTS19_02", 1, 3),
    @("1K_PFF_ADAMTS19_03", 0, 0),
    @("1K_PFF_ADAMTS19_03", 1, 6),
    @("1K_PFF_ADAMTS19_04", 0, 2),
    @("1K_PFF_ADAMTS19_04", 1, 0),
    @("1K_PFF_ADAMTS19_05", 0, 33),
    @("1K_PFF_ADAMTS19_05", 1, 8),
    @("1K_PFF_ADAMTS19_06", 0, 0),
    @("1K_PFF_ADAMTS19_06", 1, 0),
    @("1K_PFF_ADAMTS19_06", 2, 0),
    @("1K_PFF_ADAMTS19_06", 3, 0),
    @("1K_PFF_ADAMTS19_06", 4, 5),
    @("1K_PFF_ADAMTS19_06", 5, 0),
    @("1K_PFF_ADAMTS19_06", 6, 0),
    @("1K_PFF_ADAMTS19_07", 0, 9),
    @("1K_PFF_ADAMTS19_07", 1, 1),
    @("1K_PFF_ADAMTS19_08", 0, 9),
    @("1K_PFF_ADAMTS19_09", 0, 1),
    @("1K_PFF_ADAMTS19_10", 0, 0),
    @("1K_PFF_ADAMTS19_10", 1, 0),
    @("1K_PFF_ADAMTS19_10", 2, 0),
    @("1K_PFF_ADAMTS19_10", 3, 0),
    @("1K_PFF_ADAMTS19_10", 4, 0),
    @("1K_PFF_SCR_01", 0, 73),
    @("1K_PFF_SCR_01", 1, 16),
    @("1K_PFF_SCR_02", 0, 18),
    @("1K_PFF_SCR_03", 0, 0),
    @("1K_PFF_SCR_03", 1, 2),
    @("1K_PFF_SCR_03", 2, 0),
    @("1K_PFF_SCR_04", 0, 0),
    @("1K_PFF_SCR_04", 1, 1),
    @("1K_PFF_SCR_04", 2, 0),
    @("1K_PFF_SCR_05", 0, 0),
    @("1K_PFF_SCR_06", 0, 0),
    @("1K_PFF_SCR_07", 0, 0),
    @("1K_PFF_SCR_08", 0, 1),
    @("1K_PFF_SCR_08", 1, 1),
    @("1K_PFF_SCR_08", 2, 0),
    @("1K_PFF_SCR_09", 0, 0),
    @("1K_PFF_SCR_10", 0, 0),
    @("1K_PFF_TAX1BP1_01", 0, 0),
    @("1K_PFF_TAX1BP1_01", 1, 1),
    @("1K_PFF_TAX1BP1_01", 2, 1),
    @("1K_PFF_TAX1BP1_02", 0, 4),
    @("1K_PFF_TAX1BP1_02", 1, 2),
    @("1K_PFF_TAX1BP1_02", 2, 2),
    @("1K_PFF_TAX1BP1_03", 0, 0),
    @("1K_PFF_TAX1BP1_03", 1, 0),
    @("1K_PFF_TAX1BP1_04", 0, 0),
    @("1K_PFF_TAX1BP1_05", 0, 0),
    @("1K_PFF_TAX1BP1_06", 0, 1),
    @("1K_PFF_TAX1BP1_06", 1, 0),
    @("1K_PFF_TAX1BP1_06", 2, 0),
    @("1K_PFF_TAX1BP1_07", 0, 0),
    @("1K_PFF_TAX1BP1_07", 1, 0),
    @("1K_PFF_TAX1BP1_08", 0, 2),
    @("1K_PFF_TAX1BP1_09", 0, 0),
    @("1K_PFF_TAX1BP1_10", 0, 0),
    @("1K_PFF_TAX1BP1_10", 1, 0),
    @("1K_PFF_TAX1BP1_10", 2, 0),
    @("1K_PFF_TAX1BP1_10", 3, 0)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $item = $data[$i]
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
}

Write-Host "Wrote $($data.Count) data rows (rows 2..$($data.Count + 1))"
